$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new export timestamp
$ws.Name = "IClientBalance-20241125-093507-"

# Find the last used row in column G (date reference column) and update every
# data row's date value from 45618 (2024-11-22) to 45621 (2024-11-25).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 274 }

$rng = $ws.Range("G2:G$lastRow")
$rng.Value = 45621
